# Project DesignFirst save: update the computed/reported value in D10
# (Rules sheet) from 21 to 100.0, matching the upstream recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("D10").Value = 100.0
